$wb = $excel.ActiveWorkbook

# ALC row 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1586

# ALC row 13
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 533
$ws.Range("J13").Value = 533
$ws.Range("L13").Value = 533
$ws.Range("N13").Value = -871

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 9999.333000000001
$ws.Range("I86").Value = 9999
$ws.Range("J86").Value = 9999.5
$ws.Range("K86").Value = 9999
$ws.Range("L86").Value = 9999.5
$ws.Range("M86").Value = -8876
$ws.Range("N86").Value = -12245.5

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 9999.333000000001
$ws.Range("I89").Value = 9999
$ws.Range("J89").Value = 9999.5
$ws.Range("K89").Value = 49995
$ws.Range("L89").Value = 49997.5
$ws.Range("M89").Value = -44379
$ws.Range("N89").Value = -61229.5

# ALC row 93
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 700
$ws.Range("I100").Value = 700
$ws.Range("K100").Value = 700
$ws.Range("M100").Value = -159

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1850
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1850
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 5550
$ws.Range("M111").Value = ""
$ws.Range("N111").Value = -11684

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7584.879
$ws.Range("I132").Value = 8381.556
$ws.Range("J132").Value = 3999.8333
$ws.Range("K132").Value = 25144.668
$ws.Range("L132").Value = 11999.4999
$ws.Range("M132").Value = -22614.668
$ws.Range("N132").Value = -17059.4999

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5624.5
$ws.Range("I61").Value = 5714.5713
$ws.Range("K61").Value = 5714.5713
$ws.Range("M61").Value = -5502.5713

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 694.3
$ws.Range("I97").Value = 777.5714
$ws.Range("K97").Value = 777.5714
$ws.Range("M97").Value = -281.5714

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 281.84616
$ws.Range("I132").Value = 256.6
$ws.Range("J132").Value = 366
$ws.Range("K132").Value = 769.8000000000001
$ws.Range("L132").Value = 1098
$ws.Range("M132").Value = 1760.2
$ws.Range("N132").Value = -6158

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5624.5
$ws.Range("I136").Value = 5714.5713
$ws.Range("K136").Value = 17143.7139
$ws.Range("M136").Value = -14593.7139

# BSM row 14
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3137.3215
$ws.Range("I86").Value = 1274.9445
$ws.Range("J86").Value = 6489.6
$ws.Range("K86").Value = 1274.9445
$ws.Range("L86").Value = 6489.6
$ws.Range("M86").Value = -151.9445000000001
$ws.Range("N86").Value = -8735.6

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3137.3215
$ws.Range("I89").Value = 1274.9445
$ws.Range("J89").Value = 6489.6
$ws.Range("K89").Value = 6374.7225
$ws.Range("L89").Value = 32448
$ws.Range("M89").Value = -758.7224999999999
$ws.Range("N89").Value = -43680

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1734.9166
$ws.Range("I105").Value = 1762.375
$ws.Range("J105").Value = 1680
$ws.Range("K105").Value = 1762.375
$ws.Range("L105").Value = 1680
$ws.Range("M105").Value = -15.375
$ws.Range("N105").Value = -5174

# BSM row 122
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 200000
$ws.Range("J132").Value = 200000
$ws.Range("L132").Value = 200000
$ws.Range("N132").Value = -210120

# CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 15250.25
$ws.Range("I6").Value = 2000
$ws.Range("J6").Value = 19667
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 19667
$ws.Range("M6").Value = -1887
$ws.Range("N6").Value = -19893

# CRP row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("M19").Value = -330

# CRP row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 500
$ws.Range("I24").Value = 500
$ws.Range("K24").Value = 500
$ws.Range("M24").Value = -330

# CRP row 33
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1013
$ws.Range("I33").Value = 1013
$ws.Range("K33").Value = 1013
$ws.Range("M33").Value = -634

# CRP row 36
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 550
$ws.Range("I36").Value = 550
$ws.Range("K36").Value = 550
$ws.Range("M36").Value = -162

# CRP row 40
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 550
$ws.Range("I40").Value = 550
$ws.Range("K40").Value = 550
$ws.Range("M40").Value = -390

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = ""

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = ""

# CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3499.5
$ws.Range("J58").Value = 3499.5
$ws.Range("L58").Value = 10498.5
$ws.Range("N58").Value = -10754.5

# CUL row 110
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 11634
$ws.Range("I110").Value = 11634
$ws.Range("K110").Value = 34902
$ws.Range("M110").Value = -30812

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1221.6
$ws.Range("J121").Value = 2500
$ws.Range("L121").Value = 7500
$ws.Range("N121").Value = -10120

# GSM row 23
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 600
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3033.6667
$ws.Range("I80").Value = 2800.5
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 2800.5
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -1802.5
$ws.Range("N80").Value = -5496

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3033.6667
$ws.Range("I83").Value = 2800.5
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 14002.5
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -9010.5
$ws.Range("N83").Value = -27484

# GSM row 101
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 17248.75
$ws.Range("J101").Value = 17248.75
$ws.Range("L101").Value = 17248.75
$ws.Range("N101").Value = -23738.75

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5523
$ws.Range("I102").Value = 5523
$ws.Range("K102").Value = 5523
$ws.Range("M102").Value = -3901

# GSM row 109
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 69900
$ws.Range("J109").Value = 69900
$ws.Range("L109").Value = 69900
$ws.Range("N109").Value = -71980

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2880.0908
$ws.Range("I122").Value = 2668.1
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8004.299999999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -5554.299999999999
$ws.Range("N122").Value = -19900

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2166.3333
$ws.Range("I126").Value = 1999.5
$ws.Range("K126").Value = 5998.5
$ws.Range("M126").Value = -3528.5

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4260.8184
$ws.Range("I40").Value = 4260.8184
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4260.8184
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4124.8184
$ws.Range("N40").Value = ""

# WVR row 5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 20000000
$ws.Range("J5").Value = 20000000
$ws.Range("L5").Value = 20000000
$ws.Range("N5").Value = -20000224

# WVR row 75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 60000
$ws.Range("J75").Value = 60000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61872

# WVR row 78
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 60000
$ws.Range("J78").Value = 60000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -189360

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 872.58826
$ws.Range("I100").Value = 656.46155
$ws.Range("J100").Value = 1575
$ws.Range("K100").Value = 1312.9231
$ws.Range("L100").Value = 3150
$ws.Range("M100").Value = -771.9231
$ws.Range("N100").Value = -4232

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2717.4375
$ws.Range("I122").Value = 1715.5834
$ws.Range("J122").Value = 5723
$ws.Range("K122").Value = 5146.7502
$ws.Range("L122").Value = 17169
$ws.Range("M122").Value = -2696.7502
$ws.Range("N122").Value = -22069

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4940.0435
$ws.Range("I126").Value = 3314.7334
$ws.Range("K126").Value = 9944.200199999999
$ws.Range("M126").Value = -7474.200199999999
